# Apply the YX_C.xlsx update:
#  - Header A1 "Run" -> "Model_"
#  - Remove columns AL:AX (t_value_k_N .. t_value_O2_sat), keeping only
#    AK which is re-labelled "t_value_YX_C" (formerly AQ1's label)
#  - Row 3 ("Run_2" -> "Model2") and Row 4 ("Run_3" -> "Model3") get new
#    recalculated values for B and T:AK

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns AL through AX entirely; Excel will shift the dimension
# back down to A1:AK4 automatically.
$ws.Range("AL1:AX4").EntireColumn.Delete()

# --- Row 1 (headers) ---
$ws.Cells.Item(1, 1).Value = "Model_"          # A1
$ws.Cells.Item(1, 37).Value = "t_value_YX_C"   # AK1

# --- Row 3 ("Run_2" -> "Model2") ---
$ws.Cells.Item(3, 1).Value = "Model2"               # A3
$ws.Cells.Item(3, 2).Value = 0.6                    # B3

$ws.Cells.Item(3, 20).Value = 0.200767505126583     # T3
$ws.Cells.Item(3, 21).Value = 0.2046578461235455    # U3
$ws.Cells.Item(3, 22).Value = 17.92604995921369     # V3
$ws.Cells.Item(3, 23).Value = 0.4248263748452549    # W3
$ws.Cells.Item(3, 24).Value = 0.1864527495882853    # X3
$ws.Cells.Item(3, 25).Value = 7.290600700927325     # Y3
$ws.Cells.Item(3, 26).Value = 0.06422355857024478   # Z3
$ws.Cells.Item(3, 27).Value = 0.5062937016378719    # AA3
$ws.Cells.Item(3, 28).Value = 6.511507269009203     # AB3
$ws.Cells.Item(3, 29).Value = 0.1260217073942395    # AC3
$ws.Cells.Item(3, 30).Value = 0.1595211486003031    # AD3
$ws.Cells.Item(3, 31).Value = 1.448258001907745     # AE3
$ws.Cells.Item(3, 32).Value = -99.16430457356321    # AF3
$ws.Cells.Item(3, 33).Value = -97.58078563510711    # AG3
$ws.Cells.Item(3, 34).Value = 0.2453523874394228    # AH3
$ws.Cells.Item(3, 35).Value = 0.03525016248384142   # AI3
$ws.Cells.Item(3, 36).Value = 8.294103982764494     # AJ3
$ws.Cells.Item(3, 37).Value = 180.0840091539749     # AK3

# --- Row 4 ("Run_3" -> "Model3") ---
$ws.Cells.Item(4, 1).Value = "Model3"               # A4
$ws.Cells.Item(4, 2).Value = 0.6                    # B4

$ws.Cells.Item(4, 20).Value = 0.200767505126583     # T4
$ws.Cells.Item(4, 21).Value = 0.2046578461235455    # U4
$ws.Cells.Item(4, 22).Value = 17.92604995921369     # V4
$ws.Cells.Item(4, 23).Value = 0.4248263748452549    # W4
$ws.Cells.Item(4, 24).Value = 0.1864527495882853    # X4
$ws.Cells.Item(4, 25).Value = 7.290600700927325     # Y4
$ws.Cells.Item(4, 26).Value = 0.06422355857024478   # Z4
$ws.Cells.Item(4, 27).Value = 0.5062937016378719    # AA4
$ws.Cells.Item(4, 28).Value = 6.511507269009203     # AB4
$ws.Cells.Item(4, 29).Value = 0.1260217073942395    # AC4
$ws.Cells.Item(4, 30).Value = 0.1595211486003031    # AD4
$ws.Cells.Item(4, 31).Value = 1.448258001907745     # AE4
$ws.Cells.Item(4, 32).Value = -99.16430457356321    # AF4
$ws.Cells.Item(4, 33).Value = -97.58078563510711    # AG4
$ws.Cells.Item(4, 34).Value = 0.2453523874394228    # AH4
$ws.Cells.Item(4, 35).Value = 0.03525016248384142   # AI4
$ws.Cells.Item(4, 36).Value = 8.294103982764494     # AJ4
$ws.Cells.Item(4, 37).Value = 180.0840091539749     # AK4
